# Deploy the implementation guide.
# - rename the "Include from Ferlab.bio CodeS" sheet to "Include #0"
# - refresh the Metadata sheet: new Date, new Contact text, and a new
#   "Jurisdiction" row inserted right after "Contact"

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")

# Rename the Include sheet.
$wsInclude.Name = "Include #0"

# Date value (row 8) gets bumped to the new generation timestamp.
$wsMeta.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Contact value (row 10) now shows the publisher contact line.
$wsMeta.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row.
$wsMeta.Rows(11).Insert()

# Carry over the same look (border/alignment) as the surrounding rows.
$wsMeta.Range("A12:B12").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
